$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to have two header rows: row 1 held fragments of a
# badly-split "...mation (pompes)" caption plus "Hiver"/"Eté"/"Année",
# and row 2 held the matching unit labels. Remove the stray units row
# and rebuild a single, proper header row with named columns.
$ws.Rows("2:2").Delete()

$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# A1:E1 are brand-new columns (idx/idx2/Name/Date Start/Date End) and
# should use the plain default cell style.
$idCols = $ws.Range("A1:E1")
$idCols.Font.Name = "Arial"
$idCols.Font.Size = 10

# F1:K1 keep the workbook's "header" look (Arial 9 / general format)
# that was already used elsewhere in the sheet, just without a specific
# number format applied.
$units = $ws.Range("F1:K1")
$units.Font.Name = "Arial"
$units.Font.Size = 9
$units.Orientation = 0

# Match the author's cursor/selection position after the edit.
$ws.Range("A2:K2").Select()
